$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = -7.712000000000001
    "D14" = -7.542
    "D16" = -8.375
    "D21" = -8.440999999999999
    "D23" = -7.874
    "D25" = -8.228999999999999
    "D26" = -8.238
    "D29" = -7.306
    "D40" = -8.273999999999997
    "D53" = -7.443000000000001
    "D57" = -7.904000000000001
    "D59" = -8.098000000000001
    "D65" = -7.803
    "D69" = -7.221000000000001
    "D79" = -7.885
    "D83" = -8.384
    "D91" = -6.787999999999999
    "D93" = -7.531999999999999
    "D100" = -8.182
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
